# Adds "Collected Imaging" variables (MRI/PET/CAT/Other) and a
# "Polygenic Risk Scores" variable to the IHCC data dictionary
# registration template, per commit "Add collected imaging variables
# and polygenic risk scores".

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsMetadata = $wb.Worksheets.Item("Metadata")

# ---------------------------------------------------------------
# First lay out all the new rows/formatting (structure only, no new
# text yet) on both sheets.
# ---------------------------------------------------------------

# Instructions: insert a new documentation row (row 17) that explains
# the new "Collected Imaging" field, shifting the existing Terminology
# block (old rows 17-23) down to rows 18-24.
$wsInstructions.Rows.Item(17).Insert()
$wsInstructions.Range("A6:C6").Copy()
$wsInstructions.Range("A17:C17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsInstructions.Rows.Item(17).RowHeight = 15.75

# Metadata: new "Collected Imaging" section (header + MRI/PET/CAT/Other)
# after the existing "Available Datatypes" block.
$wsMetadata.Range("A23:B23").Copy()
$wsMetadata.Range("A32:B32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsMetadata.Rows.Item(32).RowHeight = $wsMetadata.Rows.Item(23).RowHeight

$wsMetadata.Range("A24:B24").Copy()
$wsMetadata.Range("A33:B36").PasteSpecial(-4122)
$excel.CutCopyMode = $false
for ($r = 33; $r -le 36; $r++) {
    $wsMetadata.Rows.Item($r).RowHeight = $wsMetadata.Rows.Item(24).RowHeight
}

# Metadata: new "Polygenic Risk Scores" field (same look as the other
# TRUE/FALSE datatype rows, e.g. row 29).
$wsMetadata.Range("A29:B29").Copy()
$wsMetadata.Range("A30:B30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsMetadata.Rows.Item(30).RowHeight = $wsMetadata.Rows.Item(29).RowHeight

# ---------------------------------------------------------------
# Now fill in the new text, in the same order the original authors
# would have: the Metadata "Collected Imaging" section first, then
# the Instructions documentation row, then "Polygenic Risk Scores".
# ---------------------------------------------------------------

$wsMetadata.Range("A32").Value = "Collected Imaging"
$wsMetadata.Range("B32").Value = "Please select TRUE or FALSE"
$wsMetadata.Range("A33").Value = "MRI"
$wsMetadata.Range("A34").Value = "PET"
$wsMetadata.Range("A35").Value = "CAT"
$wsMetadata.Range("A36").Value = "Other"

$wsInstructions.Range("A17").Value = "Collected Imaging"
$wsInstructions.Range("B17").Value = "Required"
$wsInstructions.Range("C17").Value = "TRUE/FALSE values for the types of imaging collected by the cohort."

$wsMetadata.Range("A30").Value = "Polygenic Risk Scores"

# ---------------------------------------------------------------
# Extend the TRUE/FALSE list validation that covered B24:B29 so it
# also covers the new Polygenic Risk Scores row (B30) and the new
# Collected Imaging rows (B33:B36), while leaving the blank separator
# row (31) and the new section header row (32) unvalidated - matching
# dataValidation sqref="B24:B30 B33:B36".
# ---------------------------------------------------------------
$wsMetadata.Range("B24:B36").Validation.Delete()
$wsMetadata.Range("B24:B36").Validation.Add(3, 1, 1, "TRUE, FALSE")
$wsMetadata.Range("B31:B32").Validation.Delete()

Write-Host "Added Collected Imaging (MRI/PET/CAT/Other) and Polygenic Risk Scores."
